# Auto-generated Excel COM-interop edit script
# Applies profit/formula recalculation updates across multiple worksheets
$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1407.1428
$ws.Range("I43").Value = 1350
$ws.Range("J43").Value = 1416.6666
$ws.Range("K43").Value = 1350
$ws.Range("L43").Value = 1416.6666
$ws.Range("N43").Value = -1554.6666
$ws.Range("M43").Value = -1281
$ws.Range("H62").Value = 2036.8182
$ws.Range("I62").Value = 1124.5
$ws.Range("J62").Value = 2558.1428
$ws.Range("K62").Value = 1124.5
$ws.Range("L62").Value = 2558.1428
$ws.Range("M62").Value = -500.5
$ws.Range("N62").Value = -3806.1428
$ws.Range("H63").Value = 27137.5
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 27137.5
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 27137.5
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -28385.5
$ws.Range("H65").Value = 2036.8182
$ws.Range("I65").Value = 1124.5
$ws.Range("J65").Value = 2558.1428
$ws.Range("K65").Value = 5622.5
$ws.Range("L65").Value = 12790.714
$ws.Range("M65").Value = -2502.5
$ws.Range("N65").Value = -19030.714
$ws.Range("H66").Value = 27137.5
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 27137.5
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 81412.5
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -87652.5
$ws.Range("H100").Value = 3196.4167
$ws.Range("I100").Value = 2987.889
$ws.Range("J100").Value = 3822
$ws.Range("K100").Value = 2987.889
$ws.Range("L100").Value = 3822
$ws.Range("M100").Value = -2446.889
$ws.Range("N100").Value = -4904
$ws.Range("H135").Value = 968.7895
$ws.Range("I135").Value = 994.5294
$ws.Range("J135").Value = 750
$ws.Range("K135").Value = 8950.7646
$ws.Range("L135").Value = 6750
$ws.Range("M135").Value = -6415.7646
$ws.Range("N135").Value = -11820
$ws.Range("H137").Value = 4317.6763
$ws.Range("I137").Value = 4363
$ws.Range("K137").Value = 13089
$ws.Range("M137").Value = -10539
$ws.Range("H138").Value = 3904.5762
$ws.Range("J138").Value = 7417.375
$ws.Range("L138").Value = 22252.125
$ws.Range("N138").Value = -32532.125

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 13125.667
$ws.Range("J43").Value = 13125.667
$ws.Range("L43").Value = 13125.667
$ws.Range("N43").Value = -13751.667
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H61").Value = 9488.571
$ws.Range("I61").Value = 16744.572
$ws.Range("J61").Value = 2232.5715
$ws.Range("K61").Value = 16744.572
$ws.Range("L61").Value = 2232.5715
$ws.Range("M61").Value = -16532.572
$ws.Range("N61").Value = -2656.5715
$ws.Range("H74").Value = 1078.7646
$ws.Range("I74").Value = 788.4138
$ws.Range("K74").Value = 788.4138
$ws.Range("M74").Value = 85.58619999999996
$ws.Range("H77").Value = 1078.7646
$ws.Range("I77").Value = 788.4138
$ws.Range("K77").Value = 3942.069
$ws.Range("M77").Value = 425.9309999999996
$ws.Range("H103").Value = 25076.262
$ws.Range("J103").Value = 25076.262
$ws.Range("L103").Value = 25076.262
$ws.Range("N103").Value = -27420.262
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 1646.9714
$ws.Range("I122").Value = 1374.6086
$ws.Range("K122").Value = 4123.825800000001
$ws.Range("M122").Value = -1673.825800000001
$ws.Range("H136").Value = 9488.571
$ws.Range("I136").Value = 16744.572
$ws.Range("J136").Value = 2232.5715
$ws.Range("K136").Value = 50233.716
$ws.Range("L136").Value = 6697.7145
$ws.Range("M136").Value = -47683.716
$ws.Range("N136").Value = -11797.7145

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 28800
$ws.Range("J76").Value = 28800
$ws.Range("L76").Value = 28800
$ws.Range("N76").Value = -29430
$ws.Range("H79").Value = 28800
$ws.Range("J79").Value = 28800
$ws.Range("L79").Value = 28800
$ws.Range("N79").Value = -30984
$ws.Range("H107").Value = 1372.8572
$ws.Range("I107").Value = 1351.6666
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1351.6666
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 568.3334
$ws.Range("N107").Value = -5340

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2021.3387
$ws.Range("I31").Value = 1240.8125
$ws.Range("J31").Value = 4697.4287
$ws.Range("K31").Value = 1240.8125
$ws.Range("L31").Value = 4697.4287
$ws.Range("M31").Value = -945.8125
$ws.Range("N31").Value = -5287.4287
$ws.Range("H34").Value = 2021.3387
$ws.Range("I34").Value = 1240.8125
$ws.Range("J34").Value = 4697.4287
$ws.Range("K34").Value = 1240.8125
$ws.Range("L34").Value = 4697.4287
$ws.Range("M34").Value = -1038.8125
$ws.Range("N34").Value = -5101.4287
$ws.Range("H132").Value = 1942.7826
$ws.Range("I132").Value = 1888
$ws.Range("J132").Value = 2140
$ws.Range("K132").Value = 5664
$ws.Range("L132").Value = 6420
$ws.Range("M132").Value = -3134
$ws.Range("N132").Value = -11480
$ws.Range("H135").Value = 15700
$ws.Range("J135").Value = 15700
$ws.Range("L135").Value = 15700
$ws.Range("N135").Value = -25840

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 110.85
$ws.Range("I38").Value = 59.5
$ws.Range("J38").Value = 162.2
$ws.Range("K38").Value = 178.5
$ws.Range("L38").Value = 486.6
$ws.Range("M38").Value = 168.5
$ws.Range("N38").Value = -1180.6
$ws.Range("H87").Value = 13400
$ws.Range("I87").Value = 3000
$ws.Range("K87").Value = 9000
$ws.Range("M87").Value = -7752
$ws.Range("H90").Value = 13400
$ws.Range("I90").Value = 3000
$ws.Range("K90").Value = 27000
$ws.Range("M90").Value = -20760
$ws.Range("H131").Value = 1566.6842
$ws.Range("I131").Value = 3419.8572
$ws.Range("J131").Value = 1148.2258
$ws.Range("K131").Value = 10259.5716
$ws.Range("L131").Value = 3444.6774
$ws.Range("M131").Value = -5219.571599999999
$ws.Range("N131").Value = -13524.6774
$ws.Range("H134").Value = 3217.7144
$ws.Range("I134").Value = 2498.5715
$ws.Range("J134").Value = 3936.8572
$ws.Range("K134").Value = 7495.7145
$ws.Range("L134").Value = 11810.5716
$ws.Range("M134").Value = -2425.7145
$ws.Range("N134").Value = -21950.5716

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3781.25
$ws.Range("J80").Value = 4250
$ws.Range("L80").Value = 4250
$ws.Range("N80").Value = -6246
$ws.Range("H83").Value = 3781.25
$ws.Range("J83").Value = 4250
$ws.Range("L83").Value = 21250
$ws.Range("N83").Value = -31234
$ws.Range("H122").Value = 4695.15
$ws.Range("I122").Value = 4219.375
$ws.Range("J122").Value = 5012.3335
$ws.Range("K122").Value = 12658.125
$ws.Range("L122").Value = 15037.0005
$ws.Range("M122").Value = -10208.125
$ws.Range("N122").Value = -19937.0005

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5500
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -798
$ws.Range("H113").Value = 5500
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170

